$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("all")

$ws1.Range("L3").Value = "remove STM features based on information feature strength"
$ws1.Range("L7").Value = "feature strength ~ information content by sum of reciprocals of main diagonal elements of covariance matrix"
